# Weekly update: insert a new price record as row 23, pushing the
# existing records (old rows 23-84) down by one row (new rows 24-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23 (shifts rows 23:84 -> 24:85,
# carrying over cell formatting/styles such as the date format on column D).
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new weekly data point.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44707
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100114002
$ws.Range("G23").Value = "Camote"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 18000
$ws.Range("N23").Value = "`$/caja 15 kilos granel"
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 1200
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"
